$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.710750937461853
$ws.Range("B1").Value = 1.366769671440125
$ws.Range("C1").Value = 4.01826000213623
$ws.Range("D1").Value = 5.917943477630615
$ws.Range("E1").Value = 1.823036074638367
